# The "Resources" heading paragraph used to start with a
# <w:lastRenderedPageBreak/> marker immediately followed by the
# "Resources " run. The edit splits that single paragraph into two:
#   1) the original (now empty) paragraph keeps its paragraph mark
#      formatting (Arial rFonts) but loses its run / page-break marker.
#   2) a brand-new paragraph (with the same paragraph-mark formatting)
#      carries the "Resources " run onward.
#
# Doing a Find/Replace that inserts a paragraph break ("^p") right
# before the word "Resources" reproduces exactly that split: Word
# keeps the pre-existing paragraph's identity/properties on the part
# before the break and starts a fresh paragraph for the text that
# follows, while the (rendering-only) lastRenderedPageBreak marker is
# naturally dropped since it no longer has anywhere to reattach.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Resources", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "^pResources", 2)
